$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rewrite text for existing columns, add new columns F:K ---
$ws.Range("A1").Value = "type"
$ws.Range("B1").Value = "model"
$ws.Range("C1").Value = "prec_mean"
$ws.Range("D1").Value = "rec_mean"
$ws.Range("E1").Value = "f1_mean"
$ws.Range("F1").Value = "valid_count"
$ws.Range("G1").Value = "invalid_count"
$ws.Range("H1").Value = "JSON range error"
$ws.Range("I1").Value = "JSON format error"
$ws.Range("J1").Value = "JSON key error"
$ws.Range("K1").Value = "ASSISTANT: length"

# Apply the existing bold/centered header style (from A1) to the new header cells F1:K1
$ws.Range("A1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (row 2 onward) ---
# Row 2
$ws.Range("A2").Value = "extraction"
$ws.Range("B2").Value = "gpt-4o-2024-08-06"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 7

# Row 3
$ws.Range("A3").Value = "generation"
$ws.Range("B3").Value = "Llama-2-13b-chat-hf"
$ws.Range("C3").Value = 0.1670083238018021
$ws.Range("D3").Value = 0.1575076091380439
$ws.Range("E3").Value = 0.1607123055166005
$ws.Range("F3").Value = 92
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

# Row 4
$ws.Range("A4").Value = "generation"
$ws.Range("B4").Value = "Llama-2-70b-chat-hf"
$ws.Range("C4").Value = 0.4194763268132833
$ws.Range("D4").Value = 0.3994817671111587
$ws.Range("E4").Value = 0.4026762770471262
$ws.Range("F4").Value = 92
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# Row 5
$ws.Range("A5").Value = "generation"
$ws.Range("B5").Value = "Llama-2-7b-chat-hf"
$ws.Range("C5").Value = 0.03246282702804442
$ws.Range("D5").Value = 0.03467908902691512
$ws.Range("E5").Value = 0.03082726493895909
$ws.Range("F5").Value = 92
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0

# Row 6
$ws.Range("A6").Value = "generation"
$ws.Range("B6").Value = "Llama-3-70b-chat-hf"
$ws.Range("C6").Value = 0.7344229844229845
$ws.Range("D6").Value = 0.6678618951152159
$ws.Range("E6").Value = 0.6936316051517727
$ws.Range("F6").Value = 93
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0

# Row 7
$ws.Range("A7").Value = "generation"
$ws.Range("B7").Value = "Llama-3-8b-chat-hf"
$ws.Range("C7").Value = 0.443705846125201
$ws.Range("D7").Value = 0.4268400937564359
$ws.Range("E7").Value = 0.4330634006234316
$ws.Range("F7").Value = 93
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

# Row 8
$ws.Range("A8").Value = "generation"
$ws.Range("B8").Value = "Meta-Llama-3-70B-Instruct-Lite"
$ws.Range("C8").Value = 0.7051958644700581
$ws.Range("D8").Value = 0.6080602635543711
$ws.Range("E8").Value = 0.6418570086017731
$ws.Range("F8").Value = 93
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

# Row 9
$ws.Range("A9").Value = "generation"
$ws.Range("B9").Value = "Meta-Llama-3-70B-Instruct-Turbo"
$ws.Range("C9").Value = 0.7041053629763306
$ws.Range("D9").Value = 0.6312520891505712
$ws.Range("E9").Value = 0.6588515709186034
$ws.Range("F9").Value = 93
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0

# Row 10
$ws.Range("A10").Value = "generation"
$ws.Range("B10").Value = "Meta-Llama-3-8B-Instruct-Lite"
$ws.Range("C10").Value = 0.3556728665599633
$ws.Range("D10").Value = 0.3578974772331742
$ws.Range("E10").Value = 0.3519824396247602
$ws.Range("F10").Value = 93
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0

# Row 11
$ws.Range("A11").Value = "generation"
$ws.Range("B11").Value = "Meta-Llama-3-8B-Instruct-Turbo"
$ws.Range("C11").Value = 0.3965948328851555
$ws.Range("D11").Value = 0.3829032416065107
$ws.Range("E11").Value = 0.3879828507908498
$ws.Range("F11").Value = 93
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0

# Row 12
$ws.Range("A12").Value = "generation"
$ws.Range("B12").Value = "Meta-Llama-3.1-405B-Instruct-Turbo"
$ws.Range("C12").Value = 0.8068464092657641
$ws.Range("D12").Value = 0.7637910924201248
$ws.Range("E12").Value = 0.7789730883840239
$ws.Range("F12").Value = 93
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0

# Row 13
$ws.Range("A13").Value = "generation"
$ws.Range("B13").Value = "Meta-Llama-3.1-70B-Instruct-Turbo"
$ws.Range("C13").Value = 0.5984636480604223
$ws.Range("D13").Value = 0.5631499043750617
$ws.Range("E13").Value = 0.5770086804934944
$ws.Range("F13").Value = 93
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# Row 14
$ws.Range("A14").Value = "generation"
$ws.Range("B14").Value = "Meta-Llama-3.1-8B-Instruct-Turbo"
$ws.Range("C14").Value = 0.3663480380747933
$ws.Range("D14").Value = 0.3466171710462466
$ws.Range("E14").Value = 0.354204098387242
$ws.Range("F14").Value = 93
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0

